# Edit script: insert a new data row for "Perejil" (Feria Lagunitas de Puerto Montt)
# at sheet row 243, pushing all subsequent rows down by one (the former last row,
# 358, becomes row 359). This matches the target diff, which shows the whole block
# of rows 243-358 shifting down by one row and a brand-new row 359 appearing at the
# end (containing what used to be row 358's data).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new, blank row at position 243. Excel shifts rows 243..358 down to 244..359,
# carrying their formatting and values with them, and the sheet dimension grows to R359.
$ws.Rows(243).Insert()

# Populate the freshly inserted row 243 with the new record's data.
$ws.Cells.Item(243, 1).Value2 = 4
$ws.Cells.Item(243, 2).Value2 = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(243, 3).Value2 = "Los Lagos"
$ws.Cells.Item(243, 4).Value2 = 44992
$ws.Cells.Item(243, 5).Value2 = 10
$ws.Cells.Item(243, 6).Value2 = 100112044
$ws.Cells.Item(243, 7).Value2 = "Perejil"
$ws.Cells.Item(243, 8).Value2 = "Sin especificar"
$ws.Cells.Item(243, 9).Value2 = "Primera"
$ws.Cells.Item(243, 10).Value2 = 160
$ws.Cells.Item(243, 11).Value2 = 7000
$ws.Cells.Item(243, 12).Value2 = 7000
$ws.Cells.Item(243, 13).Value2 = 7000
$ws.Cells.Item(243, 14).Value2 = "`$/docena de atados (2 kilos)"
$ws.Cells.Item(243, 15).Value2 = "Región de La Araucanía"
$ws.Cells.Item(243, 16).Value2 = 3500
$ws.Cells.Item(243, 17).Value2 = 2
$ws.Cells.Item(243, 18).Value2 = "Hortaliza"
